# Generate Report for Handoff
# Replace the old GUID-based file name (814213cb-dc5f-4fa1-ad3d-e4cce8109929)
# with the newly generated one (d7e6c88a-65dd-4bbd-999e-088005dfc57f) across
# all three sheets, refresh the handoff timestamps, and keep each
# hyperlink's target address untouched while updating its display text.

$wb = $excel.ActiveWorkbook

$oldGuid = "814213cb-dc5f-4fa1-ad3d-e4cce8109929"
$newGuid = "d7e6c88a-65dd-4bbd-999e-088005dfc57f"

$oldHash = "c811a3ab1e8fca73a6f3610460190bbd09a23151"
$newHash = "7d5373166770bb8e16e93d9122d5196a465664d4"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

$wsOverview.Range("G2").Value = "2016-12-16 10:30:40"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"

$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-12-16 10:30:26"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"

$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-12-16 10:30:40"
